$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fields")

# Normalize the "Emirate/Address/Registration Place" values from ALL CAPS to Title Case,
# and fill in the newly-added "vehicle Risk Location" value.
$ws.Range("F2").Value = "Abu Dhabi"
$ws.Range("G2").Value = "Abu Dhabi"
$ws.Range("AA2").Value = "Abu Dhabi"
$ws.Range("AL2").Value = "Abu Dhabi"

$ws.Range("AL3").Select()
